$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 118 (shifts existing rows 118-187 down to 119-188)
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new data
$ws.Range("A118").Value = 5
$ws.Range("B118").Value = "Macroferia Regional de Talca"
$ws.Range("C118").Value = "Maule"
$ws.Range("D118").Value = 44452
$ws.Range("D118").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E118").Value = 7
$ws.Range("F118").Value = 100112023
$ws.Range("G118").Value = "Brócoli"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 600
$ws.Range("L118").Value = 600
$ws.Range("M118").Value = 600
$ws.Range("N118").Value = "$/unidad"
$ws.Range("O118").Value = "Región del Maule"
$ws.Range("P118").Value = 600
$ws.Range("Q118").Value = 1
$ws.Range("R118").Value = "Hortaliza"
